$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as TEXT (never auto-converted to a number/date) while
# leaving the cell style/number-format as the default "General" afterwards.
function Set-TextValue($addr, $value) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# Row 43: adopt content previously on row 44
$ws.Range("A43").Value = 130719248
Set-TextValue "I43" "50"
$ws.Range("Q43").Value = 445790
$ws.Range("R43").Value = 7026392
$ws.Range("S43").Value = 8
$ws.Range("Z43").Value = "13:01"
$ws.Range("AB43").Value = "13:01"
$ws.Range("AC43").Value = "Rikligt på gammal levande gran (13 cm dbh, ca 175 år) i gammal granskog"
$ws.Range("AJ43").Value = "gran"
$ws.Range("AK43").Value = "Picea abies"
$ws.Range("AO43").Value = "Picea abies"

# Row 44: adopt content previously on row 45
$ws.Range("A44").Value = 130719785
$ws.Range("B44").Value = 78255
$ws.Range("E44").Value = 228579
$ws.Range("F44").Value = "Liten svartspik"
$ws.Range("G44").Value = "Chaenothecopsis nana"
$ws.Range("H44").Value = "Tibell"
$ws.Range("I44").ClearContents()
$ws.Range("J44").ClearContents()
$ws.Range("Q44").Value = 445739
$ws.Range("R44").Value = 7026240
$ws.Range("S44").Value = 10
$ws.Range("Z44").Value = "13:48"
$ws.Range("AB44").Value = "13:48"
$ws.Range("AC44").Value = "På bark på stam av levande gammal gran i gammal granskog"
$ws.Range("AJ44").ClearContents()
$ws.Range("AK44").ClearContents()
$ws.Range("AO44").ClearContents()

# Row 45: adopt content previously on row 43
$ws.Range("A45").Value = 130719742
$ws.Range("B45").Value = 79243
$ws.Range("E45").Value = 6425
$ws.Range("F45").Value = "Garnlav"
$ws.Range("G45").Value = "Alectoria sarmentosa"
$ws.Range("H45").Value = "(Ach.) Ach."
Set-TextValue "I45" "25"
$ws.Range("J45").Value = "bålar"
$ws.Range("Q45").Value = 445734
$ws.Range("R45").Value = 7026264
$ws.Range("Z45").Value = "13:43"
$ws.Range("AB45").Value = "13:43"
$ws.Range("AC45").Value = "På gammal gran i gles gammal granskog"

# Row 55: adopt content previously on row 56
$ws.Range("A55").Value = 130838768
$ws.Range("B55").Value = 79243
$ws.Range("E55").Value = 6425
$ws.Range("F55").Value = "Garnlav"
$ws.Range("G55").Value = "Alectoria sarmentosa"
$ws.Range("H55").Value = "(Ach.) Ach."
$ws.Range("Q55").Value = 445697
$ws.Range("R55").Value = 7026283
$ws.Range("Z55").Value = "11:56"
$ws.Range("AB55").Value = "11:56"
$ws.Range("AC55").Value = "På gammal gran i gammal barrblandskog"

# Row 56: adopt content previously on row 55
$ws.Range("A56").Value = 130838040
$ws.Range("B56").Value = 83228
$ws.Range("E56").Value = 1467
$ws.Range("F56").Value = "Rödbrun blekspik"
$ws.Range("G56").Value = "Sclerophora coniophaea"
$ws.Range("H56").Value = "(Norman) J.Mattsson & Middelb."
$ws.Range("Q56").Value = 445709
$ws.Range("R56").Value = 7026357
$ws.Range("Z56").Value = "11:11"
$ws.Range("AB56").Value = "11:11"
$ws.Range("AC56").Value = "På död gren i hålighet vid basen av gammal levande grov gran (42 cm dbh) i gammal granskog"

# Row 59: adopt content previously on row 61
$ws.Range("A59").Value = 130838225
$ws.Range("I59").ClearContents()
$ws.Range("J59").ClearContents()
$ws.Range("Q59").Value = 445632
$ws.Range("R59").Value = 7026388
$ws.Range("Z59").Value = "11:28"
$ws.Range("AB59").Value = "11:28"
$ws.Range("AC59").Value = "På gammal gran i gles skog nära källa"

# Row 60: adopt content previously on row 59
$ws.Range("A60").Value = 130839371
$ws.Range("Q60").Value = 445780
$ws.Range("R60").Value = 7026357
$ws.Range("S60").Value = 10
$ws.Range("Z60").Value = "12:50"
$ws.Range("AB60").Value = "12:50"
$ws.Range("AC60").Value = "På gammal gran (ca 200 år) i gles gammal granskog"

# Row 61: adopt content previously on row 62
$ws.Range("A61").Value = 130839417
$ws.Range("Q61").Value = 445781
$ws.Range("R61").Value = 7026373
$ws.Range("S61").Value = 7
$ws.Range("Z61").Value = "12:54"
$ws.Range("AB61").Value = "12:54"
$ws.Range("AC61").Value = "På gammal gran i gammal granskog"

# Row 62: adopt content previously on row 63
$ws.Range("A62").Value = 130838554
$ws.Range("Q62").Value = 445665
$ws.Range("R62").Value = 7026277
$ws.Range("S62").Value = 10
$ws.Range("Z62").Value = "11:43"
$ws.Range("AB62").Value = "11:43"
$ws.Range("AC62").Value = "På gammal levande gran i björkrik granskog"

# Row 63: adopt content previously on row 60
$ws.Range("A63").Value = 130837289
Set-TextValue "I63" "50"
$ws.Range("J63").Value = "bålar"
$ws.Range("Q63").Value = 445777
$ws.Range("R63").Value = 7026331
$ws.Range("S63").Value = 6
$ws.Range("Z63").Value = "10:34"
$ws.Range("AB63").Value = "10:34"
$ws.Range("AC63").Value = "Rikligt på gammal levande gran i gammal granskog"

# Row 68: adopt content previously on row 69
$ws.Range("A68").Value = 130839110
$ws.Range("B68").Value = 79243
$ws.Range("E68").Value = 6425
$ws.Range("F68").Value = "Garnlav"
$ws.Range("G68").Value = "Alectoria sarmentosa"
$ws.Range("H68").Value = "(Ach.) Ach."
Set-TextValue "I68" "75"
$ws.Range("J68").Value = "bålar"
$ws.Range("Q68").Value = 445730
$ws.Range("R68").Value = 7026205
$ws.Range("S68").Value = 5
$ws.Range("Z68").Value = "12:29"
$ws.Range("AB68").Value = "12:29"
$ws.Range("AC68").Value = "På gammal klen död gran i gammal granskog"

# Row 69: adopt content previously on row 68
$ws.Range("A69").Value = 130837316
$ws.Range("B69").Value = 83223
$ws.Range("E69").Value = 6440
$ws.Range("F69").Value = "Vitgrynig nållav"
$ws.Range("G69").Value = "Chaenotheca subroscida"
$ws.Range("H69").Value = "(Eitner) Zahlbr."
$ws.Range("I69").ClearContents()
$ws.Range("J69").ClearContents()
$ws.Range("Q69").Value = 445777
$ws.Range("R69").Value = 7026331
$ws.Range("S69").Value = 6
$ws.Range("Z69").Value = "10:36"
$ws.Range("AB69").Value = "10:36"
$ws.Range("AC69").Value = "Rikligt på bark på stam av levande gammal gran i gammal granskog"

# Row 76: adopt content previously on row 77
$ws.Range("A76").Value = 130837733
$ws.Range("B76").Value = 79243
$ws.Range("E76").Value = 6425
$ws.Range("F76").Value = "Garnlav"
$ws.Range("G76").Value = "Alectoria sarmentosa"
$ws.Range("H76").Value = "(Ach.) Ach."
$ws.Range("Q76").Value = 445720
$ws.Range("R76").Value = 7026343
$ws.Range("S76").Value = 10
$ws.Range("Z76").Value = "10:59"
$ws.Range("AB76").Value = "10:59"
$ws.Range("AC76").Value = "På gammal död gran i gammal granskog"

# Row 77: adopt content previously on row 78
$ws.Range("A77").Value = 130837541
$ws.Range("B77").Value = 75221
$ws.Range("D77").Value = "LC"
$ws.Range("E77").Value = 6428
$ws.Range("F77").Value = "Rostfläck"
$ws.Range("G77").Value = "Arthonia vinosa"
$ws.Range("H77").Value = "Leight."
$ws.Range("Q77").Value = 445740
$ws.Range("R77").Value = 7026322
$ws.Range("S77").Value = 8
$ws.Range("Z77").Value = "10:52"
$ws.Range("AB77").Value = "10:52"
$ws.Range("AC77").Value = "På tunna kvistar vid basen på gammal levande gran"

# Row 78: adopt content previously on row 76
$ws.Range("A78").Value = 130838833
$ws.Range("B78").Value = 89193
$ws.Range("D78").Value = "NT"
$ws.Range("E78").Value = 510
$ws.Range("F78").Value = "Doftskinn"
$ws.Range("G78").Value = "Cystostereum murrayi"
$ws.Range("H78").Value = "(Berk. & M.A.Curtis.) Pouzar"
$ws.Range("Q78").Value = 445685
$ws.Range("R78").Value = 7026259
$ws.Range("S78").Value = 4
$ws.Range("Z78").Value = "12:07"
$ws.Range("AB78").Value = "12:07"
$ws.Range("AC78").Value = "På granlåga i gammal granskog"

